$wb = $excel.ActiveWorkbook
$wsWriting = $wb.Worksheets.Item("writing")
$wsDashboard = $wb.Worksheets.Item("dashboard")

# --- Add a new day of progress (row 10) to the "writing" log ---
# Copy row 9's format down to row 10 first so the date cell (A10) keeps the
# same date-formatted style as the rest of column A instead of picking up a
# brand-new style index.
$wsWriting.Range("A9").Copy()
$wsWriting.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsWriting.Range("A10").Value = 44138
$wsWriting.Range("B10").Value = 239
$wsWriting.Range("C10").Value = 87
$wsWriting.Range("D10").Value = 528
$wsWriting.Range("E10").Value = 7638
$wsWriting.Range("F10").Value = 46
$wsWriting.Range("G10").Value = 6
$wsWriting.Range("H10").Value = 5
$wsWriting.Range("I10").Value = 5
$wsWriting.Range("J10").Formula = "=SUM(B10:I10)"
$wsWriting.Range("K10").Formula = "=J10-J9"

# --- Grow Table1 to include the new row ---
$table = $wsWriting.ListObjects.Item(1)
$table.Resize($wsWriting.Range("A1:K10"))

# --- Update the dashboard chart's series ranges to cover the new row ---
$chart = $wsDashboard.ChartObjects(1).Chart
$dailySeries = $chart.SeriesCollection(1)
$dailySeries.Formula = "=SERIES(writing!`$K`$1,writing!`$A`$2:`$A`$10,writing!`$K`$2:`$K`$10,1)"
$totalSeries = $chart.SeriesCollection(2)
$totalSeries.Formula = "=SERIES(writing!`$J`$1,writing!`$A`$2:`$A`$10,writing!`$J`$2:`$J`$10,2)"

# --- Make "writing" the active tab, selection resting on the new row ---
[void]$wsWriting.Activate()
[void]$wsWriting.Range("I10").Select()
